$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidentiality disclaimer date (A80)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-77
$ws.Range("D2").Value = 0.07449849791646478
$ws.Range("E2").Value = 0.02357723577235782
$ws.Range("D3").Value = 0.04594921988564783
$ws.Range("E3").Value = 0.02079405251502697
$ws.Range("D4").Value = 0.03640299447620893
$ws.Range("E4").Value = 0.02772849185062931
$ws.Range("D5").Value = 0.03362496365926931
$ws.Range("E5").Value = 0.01005256847735847
$ws.Range("D6").Value = 0.03128122880124044
$ws.Range("E6").Value = -0.0005855181835925372
$ws.Range("D7").Value = 0.03095910456439577
$ws.Range("E7").Value = 0.04187286949825797
$ws.Range("D8").Value = 0.02998129663727106
$ws.Range("E8").Value = 0.003684824663759523
$ws.Range("D9").Value = 0.02877507510417676
$ws.Range("E9").Value = 0.0172734314484897
$ws.Range("D10").Value = 0.02636902800658979
$ws.Range("E10").Value = 0.01399475196801214
$ws.Range("D11").Value = 0.02732706657621863
$ws.Range("E11").Value = 0.01616026213881439
$ws.Range("D12").Value = 0.0238066673127241
$ws.Range("E12").Value = -0.002487167275228908
$ws.Range("D13").Value = 0.02410960364376393
$ws.Range("E13").Value = 0.007850088630032559
$ws.Range("D14").Value = 0.02032522531253029
$ws.Range("E14").Value = 0.01129991989930179
$ws.Range("D15").Value = 0.01960965209807152
$ws.Range("E15").Value = 0.01664426346169035
$ws.Range("D16").Value = 0.02090483573989728
$ws.Range("E16").Value = 0.008272412466333146
$ws.Range("D17").Value = 0.01839984494621572
$ws.Range("E17").Value = 0.02809320159268536
$ws.Range("D18").Value = 0.0177986238976645
$ws.Range("E18").Value = 0.01330690826727055
$ws.Range("D19").Value = 0.01529435991859676
$ws.Range("E19").Value = -0.01550792821048963
$ws.Range("D20").Value = 0.01419730594049811
$ws.Range("E20").Value = -0.02193144120899371
$ws.Range("D21").Value = 0.01591849985463708
$ws.Range("E21").Value = 0.03431996249916303
$ws.Range("D22").Value = 0.01408275995736021
$ws.Range("E22").Value = 0.01381778144783929
$ws.Range("D23").Value = 0.01321523403430565
$ws.Range("E23").Value = 0.004289862724392801
$ws.Range("D24").Value = 0.01505930807248764
$ws.Range("E24").Value = 0.007876599934361739
$ws.Range("D25").Value = 0.01422889814904545
$ws.Range("E25").Value = 0.008472498433541364
$ws.Range("D26").Value = 0.01253745517976548
$ws.Range("E26").Value = 0.01450059517368252
$ws.Range("D27").Value = 0.01225797073359822
$ws.Range("E27").Value = 0.01193770258518478
$ws.Range("D28").Value = 0.0121981296637271
$ws.Range("E28").Value = 0.03082881487219225
$ws.Range("D29").Value = 0.01167264269793585
$ws.Range("E29").Value = 0.0009132420091324533
$ws.Range("D30").Value = 0.01239073553638918
$ws.Range("E30").Value = 0.03363053339590172
$ws.Range("D31").Value = 0.01278902994476209
$ws.Range("E31").Value = 0.03387133439418033
$ws.Range("D32").Value = 0.01342063184417095
$ws.Range("E32").Value = 0.02485783915515838
$ws.Range("D33").Value = 0.01105485027618955
$ws.Range("E33").Value = 0.0199211045364891
$ws.Range("D34").Value = 0.01158203314274639
$ws.Range("E34").Value = 0.009777015437392844
$ws.Range("D35").Value = 0.009619391413896696
$ws.Range("E35").Value = 0.04427653947865484
$ws.Range("D36").Value = 0.01094064347320477
$ws.Range("E36").Value = 0.005713197486193256
$ws.Range("D37").Value = 0.01072681461381917
$ws.Range("E37").Value = 0.004155750293612703
$ws.Range("D38").Value = 0.01013833704816358
$ws.Range("E38").Value = 0.02057495160943423
$ws.Range("D39").Value = 0.009224731078592887
$ws.Range("E39").Value = 0.02632629477886339
$ws.Range("D40").Value = 0.009221048551216204
$ws.Range("E40").Value = 0.00773499243315956
$ws.Range("D41").Value = 0.009347126659560034
$ws.Range("E41").Value = 0.01843384861020403
$ws.Range("D42").Value = 0.008990309138482411
$ws.Range("E42").Value = 0.01954274503885922
$ws.Range("D43").Value = 0.009711890687082083
$ws.Range("E43").Value = -0.001167466597483569
$ws.Range("D44").Value = 0.009742610718092839
$ws.Range("E44").Value = 0.03272524718005854
$ws.Range("D45").Value = 0.009223277449365248
$ws.Range("E45").Value = -0.01172576832151295
$ws.Range("D46").Value = 0.009498013373388894
$ws.Range("E46").Value = -0.000550964187327696
$ws.Range("D47").Value = 0.008833414090512647
$ws.Range("E47").Value = 0.006582411795681864
$ws.Range("D48").Value = 0.007180540750072682
$ws.Range("E48").Value = -0.001889440725545311
$ws.Range("D49").Value = 0.008273669929256711
$ws.Range("E49").Value = 0.0037715517241379
$ws.Range("D50").Value = 0.007949898245954065
$ws.Range("E50").Value = 0.01499360029255814
$ws.Range("D51").Value = 0.007894078883612753
$ws.Range("E51").Value = 0.004167740826673594
$ws.Range("D52").Value = 0.007560325612946991
$ws.Range("E52").Value = 0.009959623149394581
$ws.Range("D53").Value = 0.007098071518558
$ws.Range("E53").Value = 0.02908048330944091
$ws.Range("D54").Value = 0.007431243337532707
$ws.Range("E54").Value = 0.02241696334259213
$ws.Range("D55").Value = 0.006702684368640372
$ws.Range("E55").Value = 0.007048362611147496
$ws.Range("D56").Value = 0.006639403042930516
$ws.Range("E56").Value = 0.003911723493694508
$ws.Range("D57").Value = 0.006691346060664793
$ws.Range("E57").Value = 0.001390337156760513
$ws.Range("D58").Value = 0.006379106502568078
$ws.Range("E58").Value = -0.00533223954060702
$ws.Range("D59").Value = 0.005625787382498304
$ws.Range("E59").Value = -0.01269540502131694
$ws.Range("D60").Value = 0.006617647058823529
$ws.Range("E60").Value = 0.0174629324546951
$ws.Range("D61").Value = 0.005392479891462351
$ws.Range("E61").Value = -0.0011681193278823
$ws.Range("D62").Value = 0.005759472817133443
$ws.Range("E62").Value = 0.01332615426033112
$ws.Range("D63").Value = 0.005309816842717317
$ws.Range("E63").Value = 0.01058548693239891
$ws.Range("D64").Value = 0.004878767322414963
$ws.Range("E64").Value = 0.003654854600349644
$ws.Range("D65").Value = 0.004682236650838259
$ws.Range("E65").Value = 0.01452934845599807
$ws.Range("D66").Value = 0.004430661885841652
$ws.Range("E66").Value = 0.0003280839895012377
$ws.Range("D67").Value = 0.004353328810931292
$ws.Range("E67").Value = -0.003339121143315027
$ws.Range("D68").Value = 0.00359831379009594
$ws.Range("E68").Value = 0.01405833400662537
$ws.Range("D69").Value = 0.004121620312045741
$ws.Range("E69").Value = -0.006947873315934383
$ws.Range("D70").Value = 0.003647834092450819
$ws.Range("E70").Value = 0.0111046171829341
$ws.Range("D71").Value = 0.003179377846690571
$ws.Range("E71").Value = 0.007178127286027847
$ws.Range("D72").Value = 0.002676955131311174
$ws.Range("E72").Value = 0.01272467283291401
$ws.Range("D73").Value = 0.002613722259908906
$ws.Range("E73").Value = 0.002317303770716705
$ws.Range("D74").Value = 0.002341990502955713
$ws.Range("E74").Value = 0.01713079819588703
$ws.Range("D75").Value = 0.001907161546661498
$ws.Range("E75").Value = -0.02576219512195133
$ws.Range("D76").Value = 0.001822075782537067
$ws.Range("E76").Value = 0.04653760238272531
$ws.Range("E77").Value = 0.01385570307200301

$ws.Protect()
